# Edit script: updates the "Seguimiento Pruebas CU" table.
#  1. Row "Gestionar Plan" (table row 4): clear the green (00B050) shading on
#     all 11 cells and refresh most of the test-tracking values because the
#     use case was reimplemented.
#  2. Row "Gestionar Programa" (table row 11): a page-break now falls in the
#     middle of the row, so three cells gain a <w:lastRenderedPageBreak/>
#     marker (splitting two of the text runs in the process).
#  3. Row "Gestionar Bibliografía" (table row 12): the heading cell also
#     starts acquiring a rendered page break.

$d = $word.ActiveDocument
$t = $d.Tables(1)

$wdColorAutomatic = -16777216
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# 1. "Gestionar Plan" row -> table row 4
# ---------------------------------------------------------------------------
$planRow = 4

# Remove the green background from every cell in the row.
for ($col = 1; $col -le 11; $col++) {
    $t.Cell($planRow, $col).Shading.BackgroundPatternColor = $wdColorAutomatic
}

# Column 3: "Realizada" -> "-"
$t.Cell($planRow, 3).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 4: "26 y 27/03/2019" -> "18 y 19/08" + "/20" + "20" (three runs)
$t.Cell($planRow, 4).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr>" +
    "<w:r><w:t>18 y 19/08</w:t></w:r>" +
    "<w:r><w:t>/20</w:t></w:r>" +
    "<w:r><w:t>20</w:t></w:r></w:p>"
)

# Column 5: "18 y 19/08/2019" -> "-"
$t.Cell($planRow, 5).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 6: "SI" -> "N/A"
$t.Cell($planRow, 6).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>N/A</w:t></w:r></w:p>"
)

# Column 7: "SI" -> "-"
$t.Cell($planRow, 7).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 8: "SI" -> "-"
$t.Cell($planRow, 8).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 9: "SI" -> "-"
$t.Cell($planRow, 9).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 10: "NO" -> "-"
$t.Cell($planRow, 10).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>-</w:t></w:r></w:p>"
)

# Column 11: "Aprobado (CU cerrado - No se debe modificar)" -> "En Prueba."
$t.Cell($planRow, 11).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr><w:r><w:t>En Prueba.</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------------
# 2. "Gestionar Programa" row -> table row 11
# ---------------------------------------------------------------------------
$programaRow = 11

# Column 1: split "Gestionar Programa" into "Gestionar " + "Programa" and
# move the lastRenderedPageBreak marker to the second run.
$t.Cell($programaRow, 1).Range.InsertXML(
    "<w:p $wNs>" +
    "<w:pPr><w:jc w:val='center'/><w:rPr><w:b/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:t xml:space='preserve'>Gestionar </w:t></w:r>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Programa</w:t></w:r>" +
    "</w:p>"
)

# Column 2: "Realizada" gains a lastRenderedPageBreak marker before the text.
$t.Cell($programaRow, 2).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>Realizada</w:t></w:r></w:p>"
)

# Column 11: split the approval text around the page break.
$t.Cell($programaRow, 11).Range.InsertXML(
    "<w:p $wNs><w:pPr><w:jc w:val='center'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Aprobado (CU </w:t></w:r>" +
    "<w:r><w:lastRenderedPageBreak/><w:t>cerrado - No se debe modificar)</w:t></w:r>" +
    "</w:p>"
)

# ---------------------------------------------------------------------------
# 3. "Gestionar Bibliografía" row -> table row 12
# ---------------------------------------------------------------------------
$bibliografiaRow = 12

# Column 1: add a lastRenderedPageBreak marker before the heading text.
$t.Cell($bibliografiaRow, 1).Range.InsertXML(
    "<w:p $wNs>" +
    "<w:pPr><w:jc w:val='center'/><w:rPr><w:b/></w:rPr></w:pPr>" +
    "<w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>Gestionar Bibliografía</w:t></w:r>" +
    "</w:p>"
)
